# Final touches to static assets and templates
$wb = $excel.ActiveWorkbook

# --- Products sheet: correct stock count for "Eterno Void" (E2: 10 -> 9) ---
$products = $wb.Worksheets.Item("Products")
$products.Range("E2").Value = 9

# --- Customer_Orders sheet: append new order row (row 12) ---
$orders = $wb.Worksheets.Item("Customer_Orders")

$orders.Cells.Item(12, 1).Value = 11
$orders.Cells.Item(12, 2).Value = 3
$orders.Cells.Item(12, 3).Value = "kaizen"
$orders.Cells.Item(12, 4).Value = "boarratjabol@gmail.com"
$orders.Cells.Item(12, 5).Value = "N/A"
$orders.Cells.Item(12, 6).Value = 599
$orders.Cells.Item(12, 7).Value = 77
$orders.Cells.Item(12, 8).Value = 676
$orders.Cells.Item(12, 9).Value = "credit_card"
$orders.Cells.Item(12, 10).Value = "pending"
$orders.Cells.Item(12, 11).Value = '[{"product_id": 1, "product_name": "Eterno Void", "quantity": 1, "price": 599.0}]'
$orders.Cells.Item(12, 12).Value = "2025-11-10 02:14:24"
